$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BLS Data Series")
$ws.Range("A2:A4").EntireRow.Delete()
$ws.Activate()
$ws.Range("A2:M4").Select()
Write-Host "done"
